$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 311-312),
# pushing all existing rows 311..334 down to 313..336.
$ws.Rows("311:312").Insert()

# New row 311: Primera, Volumen 3000, precios 1300
$ws.Range("A311").Value2 = 5
$ws.Range("B311").Value = "Macroferia Regional de Talca"
$ws.Range("C311").Value = "Maule"
$ws.Range("D311").Value2 = 44783
$ws.Range("E311").Value2 = 7
$ws.Range("F311").Value2 = 100112006
$ws.Range("G311").Value = "Repollo"
$ws.Range("H311").Value = "Crespo record"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value2 = 3000
$ws.Range("K311").Value2 = 1300
$ws.Range("L311").Value2 = 1300
$ws.Range("M311").Value2 = 1300
$ws.Range("N311").Value = "$/unidad"
$ws.Range("O311").Value = "Región del Maule"
$ws.Range("P311").Value2 = 1300
$ws.Range("Q311").Value2 = 1
$ws.Range("R311").Value = "Hortaliza"

# New row 312: Segunda, Volumen 3000, precios 1000
$ws.Range("A312").Value2 = 5
$ws.Range("B312").Value = "Macroferia Regional de Talca"
$ws.Range("C312").Value = "Maule"
$ws.Range("D312").Value2 = 44783
$ws.Range("E312").Value2 = 7
$ws.Range("F312").Value2 = 100112006
$ws.Range("G312").Value = "Repollo"
$ws.Range("H312").Value = "Crespo record"
$ws.Range("I312").Value = "Segunda"
$ws.Range("J312").Value2 = 3000
$ws.Range("K312").Value2 = 1000
$ws.Range("L312").Value2 = 1000
$ws.Range("M312").Value2 = 1000
$ws.Range("N312").Value = "$/unidad"
$ws.Range("O312").Value = "Región del Maule"
$ws.Range("P312").Value2 = 1000
$ws.Range("Q312").Value2 = 1
$ws.Range("R312").Value = "Hortaliza"
